$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 232; Time = "2023-12-12 13:02:10"; Cost = 0.0012 },
    @{ Row = 233; Time = "2023-12-12 13:02:43"; Cost = 0.0018 },
    @{ Row = 234; Time = "2023-12-12 13:02:53"; Cost = 0.0004 },
    @{ Row = 235; Time = "2023-12-12 13:02:58"; Cost = 0.0004 },
    @{ Row = 236; Time = "2023-12-12 13:03:11"; Cost = 0.0006000000000000001 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Time
    $ws.Cells.Item($r.Row, 2).Value = $r.Cost
}
